$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.850.36"
$ws.Range("E2").Value = "  +2.78%  "

$ws.Range("D3").Value = "2.092.28"
$ws.Range("E3").Value = "  +2.13%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.55"
$ws.Range("E5").Value = "  +0.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("E6").Value = "  +0.85%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.29"
$ws.Range("E7").Value = "  +0.62%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0842"
$ws.Range("E10").Value = "  +0.90%  "

$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("D12").Value = "2.402.19"
$ws.Range("E12").Value = "  +2.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.01"
$ws.Range("E13").Value = "  +4.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.95"
$ws.Range("E14").Value = "  +2.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.798"
$ws.Range("E15").Value = "  +4.39%  "

$ws.Range("E16").Value = "  -0.56%  "

$ws.Range("D17").Value = "2.098.29"
$ws.Range("E17").Value = "  +2.46%  "

$ws.Range("D18").Value = "38.769.76"
$ws.Range("E18").Value = "  +2.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.60"
$ws.Range("E19").Value = "  +2.98%  "

$ws.Range("E20").Value = "  +2.18%  "

$ws.Range("E21").Value = "  +1.13%  "

$ws.Range("E22").Value = "  +2.31%  "

$ws.Range("E24").Value = "  -0.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("E25").Value = "  +2.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.06"
$ws.Range("E26").Value = "  +1.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.53"
$ws.Range("E27").Value = "  +2.15%  "

$ws.Range("E28").Value = "  +10.10%  "

$ws.Range("E29").Value = "  +13.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.18"
$ws.Range("E30").Value = "  +2.20%  "

$ws.Range("E31").Value = "  +1.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.38"
$ws.Range("E32").Value = "  +5.52%  "

$ws.Range("E33").Value = "  +2.70%  "

$ws.Range("E34").Value = "  +3.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0610"
$ws.Range("E35").Value = "  +1.27%  "

$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("E37").Value = "  +1.25%  "

$ws.Range("E38").Value = "  +2.85%  "

$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("E40").Value = "  -1.29%  "

$ws.Range("D41").Value = "1.543.02"
$ws.Range("E41").Value = "  +1.11%  "

$ws.Range("E42").Value = "  +4.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.91"
$ws.Range("E43").Value = "  +3.21%  "

$ws.Range("E44").Value = "  -0.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0923"
$ws.Range("E45").Value = "  +3.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.70"
$ws.Range("E46").Value = "  +8.47%  "

$ws.Range("E47").Value = "  +1.77%  "

$ws.Range("E48").Value = "  -1.68%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.03"
$ws.Range("E49").Value = "  +2.70%  "

$ws.Range("E50").Value = "  +0.67%  "

$ws.Range("D51").Value = "2.290.03"
$ws.Range("E51").Value = "  +2.26%  "
